$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates 44301-44303 = 2021-04-15..17)
$newRows = @(
    @{ Row = 227; A = 44301; B = 1; C = 40; D = 227.6737435255279 },
    @{ Row = 228; A = 44302; B = 3; C = 29; D = 165.0634640560077 },
    @{ Row = 229; A = 44303; B = 1; C = 24; D = 136.6042461153168 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $prevRow = $row - 1

    # Copy formatting (style) of column A from the previous row, which
    # carries the date number format / alignment / border (style index 2).
    $ws.Range("A$prevRow").Copy($ws.Range("A$row"))

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
}
